$wb = $excel.ActiveWorkbook

# PIR sheet: append rows 405-418
$ws = $wb.Worksheets.Item('PIR')
$ws.Range("A405:F418").NumberFormat = "@"
$ws.Cells.Item(405, 1).Value = '2026-02-04'
$ws.Cells.Item(405, 2).Value = '14:34:51'
$ws.Cells.Item(405, 3).Value = '14:00'
$ws.Cells.Item(405, 4).Value = 'Bathroom'
$ws.Cells.Item(405, 5).Value = 'No Motion'
$ws.Cells.Item(405, 6).Value = 'Inactive'
$ws.Cells.Item(406, 1).Value = '2026-02-04'
$ws.Cells.Item(406, 2).Value = '14:34:52'
$ws.Cells.Item(406, 3).Value = '14:00'
$ws.Cells.Item(406, 4).Value = 'Bathroom'
$ws.Cells.Item(406, 5).Value = 'Motion Detected'
$ws.Cells.Item(406, 6).Value = 'Active'
$ws.Cells.Item(407, 1).Value = '2026-02-04'
$ws.Cells.Item(407, 2).Value = '14:34:55'
$ws.Cells.Item(407, 3).Value = '14:00'
$ws.Cells.Item(407, 4).Value = 'Bathroom'
$ws.Cells.Item(407, 5).Value = 'No Motion'
$ws.Cells.Item(407, 6).Value = 'Inactive'
$ws.Cells.Item(408, 1).Value = '2026-02-04'
$ws.Cells.Item(408, 2).Value = '14:34:58'
$ws.Cells.Item(408, 3).Value = '14:00'
$ws.Cells.Item(408, 4).Value = 'Bathroom'
$ws.Cells.Item(408, 5).Value = 'Motion Detected'
$ws.Cells.Item(408, 6).Value = 'Active'
$ws.Cells.Item(409, 1).Value = '2026-02-04'
$ws.Cells.Item(409, 2).Value = '14:35:02'
$ws.Cells.Item(409, 3).Value = '14:00'
$ws.Cells.Item(409, 4).Value = 'Bathroom'
$ws.Cells.Item(409, 5).Value = 'No Motion'
$ws.Cells.Item(409, 6).Value = 'Inactive'
$ws.Cells.Item(410, 1).Value = '2026-02-04'
$ws.Cells.Item(410, 2).Value = '14:35:07'
$ws.Cells.Item(410, 3).Value = '14:00'
$ws.Cells.Item(410, 4).Value = 'Bathroom'
$ws.Cells.Item(410, 5).Value = 'No Motion'
$ws.Cells.Item(410, 6).Value = 'Inactive'
$ws.Cells.Item(411, 1).Value = '2026-02-04'
$ws.Cells.Item(411, 2).Value = '14:35:12'
$ws.Cells.Item(411, 3).Value = '14:00'
$ws.Cells.Item(411, 4).Value = 'Bathroom'
$ws.Cells.Item(411, 5).Value = 'No Motion'
$ws.Cells.Item(411, 6).Value = 'Inactive'
$ws.Cells.Item(412, 1).Value = '2026-02-04'
$ws.Cells.Item(412, 2).Value = '14:35:16'
$ws.Cells.Item(412, 3).Value = '14:00'
$ws.Cells.Item(412, 4).Value = 'Bathroom'
$ws.Cells.Item(412, 5).Value = 'Motion Detected'
$ws.Cells.Item(412, 6).Value = 'Active'
$ws.Cells.Item(413, 1).Value = '2026-02-04'
$ws.Cells.Item(413, 2).Value = '14:35:24'
$ws.Cells.Item(413, 3).Value = '14:00'
$ws.Cells.Item(413, 4).Value = 'Bathroom'
$ws.Cells.Item(413, 5).Value = 'No Motion'
$ws.Cells.Item(413, 6).Value = 'Inactive'
$ws.Cells.Item(414, 1).Value = '2026-02-04'
$ws.Cells.Item(414, 2).Value = '14:35:25'
$ws.Cells.Item(414, 3).Value = '14:00'
$ws.Cells.Item(414, 4).Value = 'Bathroom'
$ws.Cells.Item(414, 5).Value = 'Motion Detected'
$ws.Cells.Item(414, 6).Value = 'Active'
$ws.Cells.Item(415, 1).Value = '2026-02-04'
$ws.Cells.Item(415, 2).Value = '14:35:32'
$ws.Cells.Item(415, 3).Value = '14:00'
$ws.Cells.Item(415, 4).Value = 'Bathroom'
$ws.Cells.Item(415, 5).Value = 'No Motion'
$ws.Cells.Item(415, 6).Value = 'Inactive'
$ws.Cells.Item(416, 1).Value = '2026-02-04'
$ws.Cells.Item(416, 2).Value = '14:35:36'
$ws.Cells.Item(416, 3).Value = '14:00'
$ws.Cells.Item(416, 4).Value = 'Bathroom'
$ws.Cells.Item(416, 5).Value = 'Motion Detected'
$ws.Cells.Item(416, 6).Value = 'Active'
$ws.Cells.Item(417, 1).Value = '2026-02-04'
$ws.Cells.Item(417, 2).Value = '14:35:46'
$ws.Cells.Item(417, 3).Value = '14:00'
$ws.Cells.Item(417, 4).Value = 'Bathroom'
$ws.Cells.Item(417, 5).Value = 'No Motion'
$ws.Cells.Item(417, 6).Value = 'Inactive'
$ws.Cells.Item(418, 1).Value = '2026-02-04'
$ws.Cells.Item(418, 2).Value = '14:35:47'
$ws.Cells.Item(418, 3).Value = '14:00'
$ws.Cells.Item(418, 4).Value = 'Bathroom'
$ws.Cells.Item(418, 5).Value = 'Motion Detected'
$ws.Cells.Item(418, 6).Value = 'Active'

# Humidity sheet: append rows 333-340
$ws = $wb.Worksheets.Item('Humidity')
$ws.Range("A333:F340").NumberFormat = "@"
$ws.Cells.Item(333, 1).Value = '2026-02-04'
$ws.Cells.Item(333, 2).Value = '14:34:53'
$ws.Cells.Item(333, 3).Value = '14:00'
$ws.Cells.Item(333, 4).Value = 'Bathroom'
$ws.Cells.Item(333, 5).Value = '80.2%'
$ws.Cells.Item(333, 6).Value = 'Active'
$ws.Cells.Item(334, 1).Value = '2026-02-04'
$ws.Cells.Item(334, 2).Value = '14:34:56'
$ws.Cells.Item(334, 3).Value = '14:00'
$ws.Cells.Item(334, 4).Value = 'Bathroom'
$ws.Cells.Item(334, 5).Value = '79.2%'
$ws.Cells.Item(334, 6).Value = 'Active'
$ws.Cells.Item(335, 1).Value = '2026-02-04'
$ws.Cells.Item(335, 2).Value = '14:35:08'
$ws.Cells.Item(335, 3).Value = '14:00'
$ws.Cells.Item(335, 4).Value = 'Bathroom'
$ws.Cells.Item(335, 5).Value = '79.1%'
$ws.Cells.Item(335, 6).Value = 'Active'
$ws.Cells.Item(336, 1).Value = '2026-02-04'
$ws.Cells.Item(336, 2).Value = '14:35:19'
$ws.Cells.Item(336, 3).Value = '14:00'
$ws.Cells.Item(336, 4).Value = 'Bathroom'
$ws.Cells.Item(336, 5).Value = '79.1%'
$ws.Cells.Item(336, 6).Value = 'Active'
$ws.Cells.Item(337, 1).Value = '2026-02-04'
$ws.Cells.Item(337, 2).Value = '14:35:29'
$ws.Cells.Item(337, 3).Value = '14:00'
$ws.Cells.Item(337, 4).Value = 'Bathroom'
$ws.Cells.Item(337, 5).Value = '80.0%'
$ws.Cells.Item(337, 6).Value = 'Active'
$ws.Cells.Item(338, 1).Value = '2026-02-04'
$ws.Cells.Item(338, 2).Value = '14:35:39'
$ws.Cells.Item(338, 3).Value = '14:00'
$ws.Cells.Item(338, 4).Value = 'Bathroom'
$ws.Cells.Item(338, 5).Value = '78.4%'
$ws.Cells.Item(338, 6).Value = 'Active'
$ws.Cells.Item(339, 1).Value = '2026-02-04'
$ws.Cells.Item(339, 2).Value = '14:35:44'
$ws.Cells.Item(339, 3).Value = '14:00'
$ws.Cells.Item(339, 4).Value = 'Bathroom'
$ws.Cells.Item(339, 5).Value = '80.0%'
$ws.Cells.Item(339, 6).Value = 'Active'
$ws.Cells.Item(340, 1).Value = '2026-02-04'
$ws.Cells.Item(340, 2).Value = '14:35:49'
$ws.Cells.Item(340, 3).Value = '14:00'
$ws.Cells.Item(340, 4).Value = 'Bathroom'
$ws.Cells.Item(340, 5).Value = '80.0%'
$ws.Cells.Item(340, 6).Value = 'Active'

# Temperature sheet: append rows 333-340
$ws = $wb.Worksheets.Item('Temperature')
$ws.Range("A333:F340").NumberFormat = "@"
$ws.Cells.Item(333, 1).Value = '2026-02-04'
$ws.Cells.Item(333, 2).Value = '14:34:54'
$ws.Cells.Item(333, 3).Value = '14:00'
$ws.Cells.Item(333, 4).Value = 'Bathroom'
$ws.Cells.Item(333, 5).Value = '24.2C'
$ws.Cells.Item(333, 6).Value = 'Active'
$ws.Cells.Item(334, 1).Value = '2026-02-04'
$ws.Cells.Item(334, 2).Value = '14:34:57'
$ws.Cells.Item(334, 3).Value = '14:00'
$ws.Cells.Item(334, 4).Value = 'Bathroom'
$ws.Cells.Item(334, 5).Value = '24.2C'
$ws.Cells.Item(334, 6).Value = 'Active'
$ws.Cells.Item(335, 1).Value = '2026-02-04'
$ws.Cells.Item(335, 2).Value = '14:35:09'
$ws.Cells.Item(335, 3).Value = '14:00'
$ws.Cells.Item(335, 4).Value = 'Bathroom'
$ws.Cells.Item(335, 5).Value = '24.2C'
$ws.Cells.Item(335, 6).Value = 'Active'
$ws.Cells.Item(336, 1).Value = '2026-02-04'
$ws.Cells.Item(336, 2).Value = '14:35:20'
$ws.Cells.Item(336, 3).Value = '14:00'
$ws.Cells.Item(336, 4).Value = 'Bathroom'
$ws.Cells.Item(336, 5).Value = '24.2C'
$ws.Cells.Item(336, 6).Value = 'Active'
$ws.Cells.Item(337, 1).Value = '2026-02-04'
$ws.Cells.Item(337, 2).Value = '14:35:30'
$ws.Cells.Item(337, 3).Value = '14:00'
$ws.Cells.Item(337, 4).Value = 'Bathroom'
$ws.Cells.Item(337, 5).Value = '24.2C'
$ws.Cells.Item(337, 6).Value = 'Active'
$ws.Cells.Item(338, 1).Value = '2026-02-04'
$ws.Cells.Item(338, 2).Value = '14:35:40'
$ws.Cells.Item(338, 3).Value = '14:00'
$ws.Cells.Item(338, 4).Value = 'Bathroom'
$ws.Cells.Item(338, 5).Value = '24.2C'
$ws.Cells.Item(338, 6).Value = 'Active'
$ws.Cells.Item(339, 1).Value = '2026-02-04'
$ws.Cells.Item(339, 2).Value = '14:35:45'
$ws.Cells.Item(339, 3).Value = '14:00'
$ws.Cells.Item(339, 4).Value = 'Bathroom'
$ws.Cells.Item(339, 5).Value = '24.2C'
$ws.Cells.Item(339, 6).Value = 'Active'
$ws.Cells.Item(340, 1).Value = '2026-02-04'
$ws.Cells.Item(340, 2).Value = '14:35:50'
$ws.Cells.Item(340, 3).Value = '14:00'
$ws.Cells.Item(340, 4).Value = 'Bathroom'
$ws.Cells.Item(340, 5).Value = '24.2C'
$ws.Cells.Item(340, 6).Value = 'Active'

